$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before column D, shifting the old D/E/F columns to
#    E/F/G. This makes room for the new "Lifebase Distance accumulated" data
#    while the old D+ Difference values (old D) and the duration-in-seconds
#    values (old F) land exactly where the target layout needs them (E & G).
# ---------------------------------------------------------------------------
$ws.Range("D1:D14").Insert(-4161)  # xlShiftToRight

# ---------------------------------------------------------------------------
# 2. Re-label the header row (row 1) to the new column headers.
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Lifebase"
$ws.Range("B1").Value = "Stage"
$ws.Range("C1").Value = "Lifebase Distance (km)"
$ws.Range("D1").Value = "Lifebase Accumulated Distance Elevation (m)"
$ws.Range("E1").Value = "Lifebase Elevation Gain (m)"
$ws.Range("F1").Value = "Lifebase Accumulated Elevation (m)"
$ws.Range("G1").Value = "Lifebase Duration_seconds"

# ---------------------------------------------------------------------------
# 3. Re-label column B (the "Stage" column) with the new stage / checkpoint
#    time labels.
# ---------------------------------------------------------------------------
$ws.Range("B2").Value  = "Stage1 Time"
$ws.Range("B3").Value  = "Time Spent in Valgrisenche OUT"
$ws.Range("B4").Value  = "Stage2 Time"
$ws.Range("B5").Value  = "Time Spent in Cogne OUT"
$ws.Range("B6").Value  = "Stage3 Time"
$ws.Range("B7").Value  = "Time Spent in Donnas OUT"
$ws.Range("B8").Value  = "Stage4 Time"
$ws.Range("B9").Value  = "Time Spent in Gressoney OUT"
$ws.Range("B10").Value = "Stage5 Time"
$ws.Range("B11").Value = "Time Spent in Valtournenche OUT"
$ws.Range("B12").Value = "Stage6 Time"
$ws.Range("B13").Value = "Time Spent in Ollomont OUT"
$ws.Range("B14").Value = "Stage7 Time"

# ---------------------------------------------------------------------------
# 4. Fill the new column D with the cumulative ("accumulated") lifebase
#    distance, running down column C.
# ---------------------------------------------------------------------------
$ws.Range("D2").Value  = 48.55
$ws.Range("D3").Value  = 48.55
$ws.Range("D4").Value  = 104
$ws.Range("D5").Value  = 104
$ws.Range("D6").Value  = 149.77
$ws.Range("D7").Value  = 149.77
$ws.Range("D8").Value  = 204
$ws.Range("D9").Value  = 204
$ws.Range("D10").Value = 237.62
$ws.Range("D11").Value = 237.62
$ws.Range("D12").Value = 285.66
$ws.Range("D13").Value = 285.66
$ws.Range("D14").Value = 335.32

# ---------------------------------------------------------------------------
# 5. Column E already holds the old "D+ Difference (m)" values (shifted over
#    by the insert in step 1) - those numbers are unchanged by this edit.
# ---------------------------------------------------------------------------

# ---------------------------------------------------------------------------
# 6. Column F currently holds the old "Stage Duration" fraction-of-day values
#    (shifted over from old column E, including their number format style).
#    Replace them with the new cumulative ("accumulated") elevation values
#    and drop the old number-format style so the cells fall back to General.
# ---------------------------------------------------------------------------
$ws.Range("F2:F14").ClearFormats()

$ws.Range("F2").Value  = 4339
$ws.Range("F3").Value  = 4339
$ws.Range("F4").Value  = 9282
$ws.Range("F5").Value  = 9282
$ws.Range("F6").Value  = 12050
$ws.Range("F7").Value  = 12050
$ws.Range("F8").Value  = 17983
$ws.Range("F9").Value  = 17983
$ws.Range("F10").Value = 21077
$ws.Range("F11").Value = 21077
$ws.Range("F12").Value = 25702
$ws.Range("F13").Value = 25702
$ws.Range("F14").Value = 29608

# ---------------------------------------------------------------------------
# 7. Column G already holds the old "Stage Duration_in_seconds" values
#    (shifted over from old column F) - unchanged by this edit.
# ---------------------------------------------------------------------------
